{"js": "// The supplemental table reports, for each variable/state row, a Math\n// column and a Reading/Language-Arts column of the form\n//   \"<beta> (<lo>, <hi>), <probability>%\"\n// The posterior probability was originally rounded to whole percent (or\n// shown as the clipped \"&gt;99.9%\"); this edit corrects it to show one\n// more decimal place (e.g. \"50%\" -> \"50.0%\", \"&gt;99.9%\" -> \"99.5%\").\n//\n// Row/column indices below are 0-based into the single table in the\n// document body: column 1 = Math, column 2 = Reading/Language Arts.\nconst changes = [\n  { row: 2, col: 1, label: \"National\", before: \"50%\", after: \"50.0%\" },\n  { row: 2, col: 2, label: \"National\", before: \"50%\", after: \"50.0%\" },\n  { row: 3, col: 1, label: \"Alabama\", before: \"93%\", after: \"93.0%\" },\n  { row: 3, col: 2, label: \"Alabama\", before: \"70%\", after: \"70.1%\" },\n  { row: 4, col: 2, label: \"Florida\", before: \"95%\", after: \"95.5%\" },\n  { row: 5, col: 1, label: \"Georgia\", before: \"90%\", after: \"90.3%\" },\n  { row: 5, col: 2, label: \"Georgia\", before: \"56%\", after: \"56.5%\" },\n  { row: 6, col: 1, label: \"Louisiana\", before: \"54%\", after: \"54.5%\" },\n  { row: 6, col: 2, label: \"Louisiana\", before: \"50%\", after: \"52.2%\" },\n  { row: 7, col: 1, label: \"New Jersey\", before: \"89%\", after: \"88.8%\" },\n  { row: 7, col: 2, label: \"New Jersey\", before: \"64%\", after: \"64.0%\" },\n  { row: 8, col: 1, label: \"North Carolina\", before: \">99.9%\", after: \"99.5%\" },\n  { row: 8, col: 2, label: \"North Carolina\", before: \"78%\", after: \"78.2%\" },\n  { row: 9, col: 1, label: \"South Carolina\", before: \"97%\", after: \"97.2%\" },\n  { row: 9, col: 2, label: \"South Carolina\", before: \"60%\", after: \"59.8%\" },\n  { row: 10, col: 1, label: \"Texas\", before: \"95%\", after: \"94.9%\" },\n  { row: 10, col: 2, label: \"Texas\", before: \">99.9%\", after: \"99.9%\" },\n  { row: 12, col: 1, label: \"American Indian/Alaska Native\", before: \">99.9%\", after: \"99.8%\" },\n  { row: 16, col: 2, label: \"Percent Grade-Cohort Receiving Free Lunch\", before: \">99.9%\", after: \"99.9%\" },\n  { row: 17, col: 1, label: \"Percent Grade-Cohort Economically Disadvantaged\", before: \"94%\", after: \"93.9%\" },\n  { row: 18, col: 1, label: \"Percent County English Language Learners\", before: \"99%\", after: \"99.4%\" },\n  { row: 19, col: 1, label: \"Percent County Urban Schools\", before: \"95%\", after: \"95.2%\" },\n  { row: 19, col: 2, label: \"Percent County Urban Schools\", before: \">99.9%\", after: \"99.9%\" },\n  { row: 20, col: 1, label: \"Percent County Special Education Students\", before: \"96%\", after: \"95.9%\" },\n  { row: 22, col: 2, label: \"County Poverty Rate\", before: \"64%\", after: \"64.4%\" },\n  { row: 23, col: 1, label: \"Percent County Single Mother Households\", before: \"52%\", after: \"52.2%\" },\n  { row: 23, col: 2, label: \"Percent County Single Mother Households\", before: \"72%\", after: \"72.5%\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length < 1) {\n  throw new Error(\"expected at least one table in the document body\");\n}\nconst table = tables.items[0];\n\n// Verify row labels line up with what we expect before mutating anything,\n// so we fail loudly instead of silently editing the wrong row.\nconst labelCells = changes.map((ch) => table.getCell(ch.row, 0).body);\nlabelCells.forEach((b) => b.load(\"text\"));\nawait context.sync();\nchanges.forEach((ch, i) => {\n  const actual = labelCells[i].text.trim();\n  if (actual !== ch.label) {\n    throw new Error(\n      `row ${ch.row} label mismatch: expected \"${ch.label}\", found \"${actual}\"`\n    );\n  }\n});\n\nfor (const ch of changes) {\n  const cell = table.getCell(ch.row, ch.col);\n  const range = cell.body.getRange();\n  const results = range.search(ch.before, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `expected exactly 1 match for \"${ch.before}\" in row ${ch.row} ` +\n        `(\"${ch.label}\"), col ${ch.col}, found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(ch.after, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The supplemental table reports, for each variable/state row, a Math\n# column and a Reading/Language-Arts column of the form\n#   \"<beta> (<lo>, <hi>), <probability>%\"\n# The posterior probability was originally rounded to whole percent (or\n# shown as the clipped \">99.9%\"); this edit corrects it to show one more\n# decimal place (e.g. \"50%\" -> \"50.0%\", \">99.9%\" -> \"99.5%\").\n#\n# Word COM table Cell(row, col) is 1-based, and column 1 is the row label,\n# column 2 is Math, column 3 is Reading/Language Arts.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Replace-InCell($table, $row, $col, $expectedLabel, $before, $after) {\n  $labelCell = $table.Cell($row, 1)\n  # Cell().Range.Text carries a trailing cell-mark (chr 7) plus CR/LF; strip\n  # those before comparing against the plain-text expected label.\n  $label = $labelCell.Range.Text.TrimEnd([char]7, [char]13, [char]10)\n  if ($label -ne $expectedLabel) {\n    throw \"row $row label mismatch: expected '$expectedLabel', found '$label'\"\n  }\n\n  $cell = $table.Cell($row, $col)\n  $rng = $cell.Range\n  $txt = $rng.Text\n  $idx = $txt.IndexOf($before)\n  if ($idx -lt 0) {\n    throw \"substring '$before' not found in cell ($row,$col); cell text='$txt'\"\n  }\n\n  # Build a tightly-bounded sub-range over just the target substring and\n  # assign its .Text directly; Find/Replace scoped to a sub-range is not\n  # reliable here, but a direct bounded-range text assignment is, and it\n  # preserves the surrounding runs (and this run's own formatting).\n  $subStart = $rng.Start + $idx\n  $subEnd = $subStart + $before.Length\n  $subRng = $d.Range($subStart, $subEnd)\n  $subRng.Text = $after\n}\n\nReplace-InCell $t 3 2 \"National\" \"50%\" \"50.0%\"\nReplace-InCell $t 3 3 \"National\" \"50%\" \"50.0%\"\nReplace-InCell $t 4 2 \"Alabama\" \"93%\" \"93.0%\"\nReplace-InCell $t 4 3 \"Alabama\" \"70%\" \"70.1%\"\nReplace-InCell $t 5 3 \"Florida\" \"95%\" \"95.5%\"\nReplace-InCell $t 6 2 \"Georgia\" \"90%\" \"90.3%\"\nReplace-InCell $t 6 3 \"Georgia\" \"56%\" \"56.5%\"\nReplace-InCell $t 7 2 \"Louisiana\" \"54%\" \"54.5%\"\nReplace-InCell $t 7 3 \"Louisiana\" \"50%\" \"52.2%\"\nReplace-InCell $t 8 2 \"New Jersey\" \"89%\" \"88.8%\"\nReplace-InCell $t 8 3 \"New Jersey\" \"64%\" \"64.0%\"\nReplace-InCell $t 9 2 \"North Carolina\" \">99.9%\" \"99.5%\"\nReplace-InCell $t 9 3 \"North Carolina\" \"78%\" \"78.2%\"\nReplace-InCell $t 10 2 \"South Carolina\" \"97%\" \"97.2%\"\nReplace-InCell $t 10 3 \"South Carolina\" \"60%\" \"59.8%\"\nReplace-InCell $t 11 2 \"Texas\" \"95%\" \"94.9%\"\nReplace-InCell $t 11 3 \"Texas\" \">99.9%\" \"99.9%\"\nReplace-InCell $t 13 2 \"American Indian/Alaska Native\" \">99.9%\" \"99.8%\"\nReplace-InCell $t 17 3 \"Percent Grade-Cohort Receiving Free Lunch\" \">99.9%\" \"99.9%\"\nReplace-InCell $t 18 2 \"Percent Grade-Cohort Economically Disadvantaged\" \"94%\" \"93.9%\"\nReplace-InCell $t 19 2 \"Percent County English Language Learners\" \"99%\" \"99.4%\"\nReplace-InCell $t 20 2 \"Percent County Urban Schools\" \"95%\" \"95.2%\"\nReplace-InCell $t 20 3 \"Percent County Urban Schools\" \">99.9%\" \"99.9%\"\nReplace-InCell $t 21 2 \"Percent County Special Education Students\" \"96%\" \"95.9%\"\nReplace-InCell $t 23 3 \"County Poverty Rate\" \"64%\" \"64.4%\"\nReplace-InCell $t 24 2 \"Percent County Single Mother Households\" \"52%\" \"52.2%\"\nReplace-InCell $t 24 3 \"Percent County Single Mother Households\" \"72%\" \"72.5%\"\n"}
